$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the redundant numeric "ID Calidad" values (MLM131441 / MLM456634)
# that were duplicated from column A into column B for the two product
# rows. No other cell content moves - only B2 and B3 are cleared.

$ws.Range("B2:B3").ClearContents()
